$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Progress value for the BNT323-01 trial (row 6) from 100 to 88
$ws.Range("B6").Value = 88
